# [PHOENIX-5854] UI trade license changes
#
# - Update the trade commencement date on the "tradeDetails" sheet (H2)
#   from 01/03/2017 to 01/04/2017.
# - Move the active selection on "tradeLocationDetails" from C8 to B10.
# - Make "tradeDetails" the active sheet (instead of "tradeLocationDetails"),
#   with its selection moved from G9 to H7.

$wb = $excel.ActiveWorkbook

$wsLocation = $wb.Worksheets.Item("tradeLocationDetails")
$wsDetails  = $wb.Worksheets.Item("tradeDetails")

# Update the trade commencement date value.
$wsDetails.Range("H2").Value = "01/04/2017"

# Update the (now inactive) selection left on tradeLocationDetails.
$wsLocation.Range("B10").Select()

# Activate tradeDetails and move its selection.
$wsDetails.Activate()
$wsDetails.Range("H7").Select()
